# Update the "Repo Link" entry for the third person (row 3) with the
# author's new repository URL, make it a hyperlink like the other
# email/link cells, give it the same "Hyperlink" cell style used by the
# sibling cells, and move the active selection to that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://github.com/EngNoOne/code_stego.git"

# 1. Correct the repo link text in C3.
$ws.Range("C3").Value = $newUrl

# 2. Turn C3 into a hyperlink pointing at the new URL (adds a new
#    relationship + <hyperlink> entry, matching the pattern already used
#    by B2/B3/C2).
$ws.Hyperlinks.Add($ws.Range("C3"), $newUrl)

# 3. Re-apply the same visual style the other hyperlink cells use so C3
#    matches B2/B3/C2 instead of picking up a freshly minted style index.
$ws.Range("C3").Style = $ws.Range("C2").Style

# 4. Move/save the current selection to C3 (previously C5).
$ws.Range("C3").Select()
